# The commit swaps the contents of ppt/theme/theme1.xml (the "Integral"
# colour theme used by the slide master) and ppt/theme/theme2.xml (the
# default "Office Theme" colours, previously only used by the notes
# master) - i.e. the slide master ends up with the plain Office Theme
# palette instead of the custom Integral palette.
#
# The PowerPoint object model only exposes the document's theme colours
# through ThemeColorScheme (12 slots: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) reachable from a Slide/SlideRange/CustomLayout, and
# that maps onto the single slide-master theme part (theme1.xml) used
# throughout the deck. Re-point every one of those 12 slots at the
# standard "Office Theme" RGB values so the master's theme becomes the
# Office Theme colour scheme, matching the target edit.

function ToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Office Theme colour scheme, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $tcs.Colors($i).RGB = ToComRgb $officeTheme[$i - 1]
}
